$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 349
$ws.Cells.Item(19, 10).Value = 342.625
$ws.Cells.Item(19, 12).Value = 342.625
$ws.Cells.Item(19, 14).Value = -692.625
$ws.Cells.Item(43, 8).Value = 9663
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 13).ClearContents()
$ws.Cells.Item(104, 8).Value = 348
$ws.Cells.Item(104, 9).Value = 348
$ws.Cells.Item(104, 11).Value = 1044
$ws.Cells.Item(104, 13).Value = 703
$ws.Cells.Item(137, 8).Value = 8106.72
$ws.Cells.Item(137, 9).Value = 12320.538
$ws.Cells.Item(137, 11).Value = 36961.614
$ws.Cells.Item(137, 13).Value = -34411.614
$ws.Cells.Item(138, 8).Value = 3583.4043
$ws.Cells.Item(138, 9).Value = 2230.7307
$ws.Cells.Item(138, 10).Value = 4100.603
$ws.Cells.Item(138, 11).Value = 6692.1921
$ws.Cells.Item(138, 12).Value = 12301.809
$ws.Cells.Item(138, 13).Value = -1552.1921
$ws.Cells.Item(138, 14).Value = -22581.809
$ws.Cells.Item(141, 8).Value = 3099.3
$ws.Cells.Item(141, 9).Value = 3165.5
$ws.Cells.Item(141, 11).Value = 9496.5
$ws.Cells.Item(141, 13).Value = -4316.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3668.3242
$ws.Cells.Item(32, 9).Value = 3840.3333
$ws.Cells.Item(32, 11).Value = 3840.3333
$ws.Cells.Item(32, 13).Value = -3553.3333
$ws.Cells.Item(61, 8).Value = 4458
$ws.Cells.Item(61, 9).Value = 4547
$ws.Cells.Item(61, 11).Value = 4547
$ws.Cells.Item(61, 13).Value = -4335
$ws.Cells.Item(74, 8).Value = 2136.84
$ws.Cells.Item(74, 9).Value = 1746.409
$ws.Cells.Item(74, 10).Value = 5000
$ws.Cells.Item(74, 11).Value = 1746.409
$ws.Cells.Item(74, 12).Value = 5000
$ws.Cells.Item(74, 13).Value = -872.4090000000001
$ws.Cells.Item(74, 14).Value = -6748
$ws.Cells.Item(77, 8).Value = 2136.84
$ws.Cells.Item(77, 9).Value = 1746.409
$ws.Cells.Item(77, 10).Value = 5000
$ws.Cells.Item(77, 11).Value = 8732.045
$ws.Cells.Item(77, 12).Value = 25000
$ws.Cells.Item(77, 13).Value = -4364.045
$ws.Cells.Item(77, 14).Value = -33736
$ws.Cells.Item(110, 8).Value = 362.45456
$ws.Cells.Item(110, 9).Value = 383.3
$ws.Cells.Item(110, 11).Value = 383.3
$ws.Cells.Item(110, 13).Value = 1661.7
$ws.Cells.Item(122, 8).Value = 3917.8
$ws.Cells.Item(122, 9).Value = 4446.3335
$ws.Cells.Item(122, 11).Value = 13339.0005
$ws.Cells.Item(122, 13).Value = -10889.0005
$ws.Cells.Item(132, 8).Value = 2627.4866
$ws.Cells.Item(132, 9).Value = 2450.0881
$ws.Cells.Item(132, 10).Value = 4638
$ws.Cells.Item(132, 11).Value = 7350.2643
$ws.Cells.Item(132, 12).Value = 13914
$ws.Cells.Item(132, 13).Value = -4820.2643
$ws.Cells.Item(132, 14).Value = -18974
$ws.Cells.Item(136, 8).Value = 4458
$ws.Cells.Item(136, 9).Value = 4547
$ws.Cells.Item(136, 11).Value = 13641
$ws.Cells.Item(136, 13).Value = -11091
$ws.Cells.Item(137, 8).Value = 149999
$ws.Cells.Item(137, 10).Value = 149999
$ws.Cells.Item(137, 12).Value = 149999
$ws.Cells.Item(137, 14).Value = -160199
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 4372
$ws.Cells.Item(134, 9).Value = 4543.7
$ws.Cells.Item(134, 11).Value = 13631.1
$ws.Cells.Item(134, 13).Value = -11096.1
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1454
$ws.Cells.Item(22, 9).Value = 969.3333
$ws.Cells.Item(22, 11).Value = 969.3333
$ws.Cells.Item(22, 13).Value = -619.3333
$ws.Cells.Item(31, 8).Value = 3173.5293
$ws.Cells.Item(31, 9).Value = 1578.95
$ws.Cells.Item(31, 11).Value = 1578.95
$ws.Cells.Item(31, 13).Value = -1283.95
$ws.Cells.Item(34, 8).Value = 3173.5293
$ws.Cells.Item(34, 9).Value = 1578.95
$ws.Cells.Item(34, 11).Value = 1578.95
$ws.Cells.Item(34, 13).Value = -1376.95
$ws.Cells.Item(82, 8).Value = 181700
$ws.Cells.Item(82, 10).Value = 182550
$ws.Cells.Item(82, 12).Value = 182550
$ws.Cells.Item(82, 14).Value = -183272
$ws.Cells.Item(85, 8).Value = 181700
$ws.Cells.Item(85, 10).Value = 182550
$ws.Cells.Item(85, 12).Value = 182550
$ws.Cells.Item(85, 14).Value = -185046
$ws.Cells.Item(134, 8).Value = 1943.5555
$ws.Cells.Item(134, 9).Value = 1778.48
$ws.Cells.Item(134, 10).Value = 4007
$ws.Cells.Item(134, 11).Value = 5335.440000000001
$ws.Cells.Item(134, 12).Value = 12021
$ws.Cells.Item(134, 13).Value = -2800.440000000001
$ws.Cells.Item(134, 14).Value = -17091
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(47, 8).Value = 496.66666
$ws.Cells.Item(47, 9).Value = 570
$ws.Cells.Item(47, 10).Value = 350
$ws.Cells.Item(47, 11).Value = 1710
$ws.Cells.Item(47, 12).Value = 1050
$ws.Cells.Item(47, 13).Value = -1279
$ws.Cells.Item(47, 14).Value = -1912
$ws.Cells.Item(57, 8).Value = 4888.4614
$ws.Cells.Item(57, 9).Value = 50
$ws.Cells.Item(57, 10).Value = 5291.6665
$ws.Cells.Item(57, 11).Value = 150
$ws.Cells.Item(57, 12).Value = 15874.9995
$ws.Cells.Item(57, 13).Value = 409
$ws.Cells.Item(57, 14).Value = -16992.9995
$ws.Cells.Item(69, 8).Value = 16671416
$ws.Cells.Item(69, 9).Value = 5249
$ws.Cells.Item(69, 10).Value = 25004500
$ws.Cells.Item(69, 11).Value = 15747
$ws.Cells.Item(69, 12).Value = 75013500
$ws.Cells.Item(69, 13).Value = -14936
$ws.Cells.Item(69, 14).Value = -75015122
$ws.Cells.Item(72, 8).Value = 16671416
$ws.Cells.Item(72, 9).Value = 5249
$ws.Cells.Item(72, 10).Value = 25004500
$ws.Cells.Item(72, 11).Value = 47241
$ws.Cells.Item(72, 12).Value = 225040500
$ws.Cells.Item(72, 13).Value = -43185
$ws.Cells.Item(72, 14).Value = -225048612
$ws.Cells.Item(113, 8).Value = 1177.3
$ws.Cells.Item(113, 9).Value = 1356.4286
$ws.Cells.Item(113, 10).Value = 1080.8462
$ws.Cells.Item(113, 11).Value = 4069.2858
$ws.Cells.Item(113, 12).Value = 3242.5386
$ws.Cells.Item(113, 13).Value = -1899.2858
$ws.Cells.Item(113, 14).Value = -7582.5386
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 10718.96
$ws.Cells.Item(80, 9).Value = 12698.267
$ws.Cells.Item(80, 11).Value = 12698.267
$ws.Cells.Item(80, 13).Value = -11700.267
$ws.Cells.Item(83, 8).Value = 10718.96
$ws.Cells.Item(83, 9).Value = 12698.267
$ws.Cells.Item(83, 11).Value = 63491.335
$ws.Cells.Item(83, 13).Value = -58499.335
$ws.Cells.Item(93, 8).Value = 33332.668
$ws.Cells.Item(93, 10).Value = 33332.668
$ws.Cells.Item(93, 12).Value = 33332.668
$ws.Cells.Item(93, 14).Value = -37076.668
$ws.Cells.Item(122, 8).Value = 6085.9644
$ws.Cells.Item(122, 9).Value = 5919.7827
$ws.Cells.Item(122, 10).Value = 6850.4
$ws.Cells.Item(122, 11).Value = 17759.3481
$ws.Cells.Item(122, 12).Value = 20551.2
$ws.Cells.Item(122, 13).Value = -15309.3481
$ws.Cells.Item(122, 14).Value = -25451.2
$ws.Cells.Item(132, 8).Value = 2950
$ws.Cells.Item(132, 9).Value = 2873
$ws.Cells.Item(132, 10).Value = 3209.875
$ws.Cells.Item(132, 11).Value = 8619
$ws.Cells.Item(132, 12).Value = 9629.625
$ws.Cells.Item(132, 13).Value = -6089
$ws.Cells.Item(132, 14).Value = -14689.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 3912.7727
$ws.Cells.Item(16, 9).Value = 3605.8667
$ws.Cells.Item(16, 11).Value = 3605.8667
$ws.Cells.Item(16, 13).Value = -3435.8667
$ws.Cells.Item(75, 8).Value = 50000
$ws.Cells.Item(75, 10).Value = 50000
$ws.Cells.Item(75, 12).Value = 50000
$ws.Cells.Item(75, 14).Value = -51872
$ws.Cells.Item(78, 8).Value = 50000
$ws.Cells.Item(78, 10).Value = 50000
$ws.Cells.Item(78, 12).Value = 150000
$ws.Cells.Item(78, 14).Value = -159360
$ws.Cells.Item(132, 8).Value = 3902.697
$ws.Cells.Item(132, 9).Value = 3681.1904
$ws.Cells.Item(132, 10).Value = 4290.3335
$ws.Cells.Item(132, 11).Value = 11043.5712
$ws.Cells.Item(132, 12).Value = 12871.0005
$ws.Cells.Item(132, 13).Value = -8513.5712
$ws.Cells.Item(132, 14).Value = -17931.0005
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1350
$ws.Cells.Item(81, 9).Value = 1700
$ws.Cells.Item(81, 10).Value = 1175
$ws.Cells.Item(81, 11).Value = 3400
$ws.Cells.Item(81, 12).Value = 2350
$ws.Cells.Item(81, 13).Value = -2339
$ws.Cells.Item(81, 14).Value = -4472
$ws.Cells.Item(84, 8).Value = 1350
$ws.Cells.Item(84, 9).Value = 1700
$ws.Cells.Item(84, 10).Value = 1175
$ws.Cells.Item(84, 11).Value = 17000
$ws.Cells.Item(84, 12).Value = 11750
$ws.Cells.Item(84, 13).Value = -11696
$ws.Cells.Item(84, 14).Value = -22358
$ws.Cells.Item(96, 8).Value = 4171.625
$ws.Cells.Item(96, 9).Value = 5044.5713
$ws.Cells.Item(96, 11).Value = 5044.5713
$ws.Cells.Item(96, 13).Value = -3671.5713
$ws.Cells.Item(108, 8).Value = 106284
$ws.Cells.Item(108, 10).Value = 106284
$ws.Cells.Item(108, 12).Value = 106284
$ws.Cells.Item(108, 14).Value = -113964
$ws.Cells.Item(126, 8).Value = 2250.04
$ws.Cells.Item(126, 9).Value = 2006.7059
$ws.Cells.Item(126, 11).Value = 6020.1177
$ws.Cells.Item(126, 13).Value = -3550.1177
$ws.Cells.Item(132, 8).Value = 3676.7666
$ws.Cells.Item(132, 9).Value = 2804.44
$ws.Cells.Item(132, 10).Value = 8038.4
$ws.Cells.Item(132, 11).Value = 8413.32
$ws.Cells.Item(132, 12).Value = 24115.2
$ws.Cells.Item(132, 13).Value = -5883.32
$ws.Cells.Item(132, 14).Value = -29175.2
